$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1914414414414414
$ws.Range("C2").Value = 0.5945945945945946
$ws.Range("J2").Value = 0.01576576576576576
$ws.Range("P2").Value = 0.1306306306306306
$ws.Range("S2").Value = 0.06756756756756757
$ws.Range("C3").Value = 0.01805054151624549
$ws.Range("J3").Value = 0.03971119133574007
$ws.Range("P3").Value = 0.7364620938628159
$ws.Range("S3").Value = 0.2057761732851986
$ws.Range("B6").Value = 0.07037037037037037
$ws.Range("D6").Value = 0.01111111111111111
$ws.Range("F6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.3148148148148148
$ws.Range("O6").Value = 0.01851851851851852
$ws.Range("Q6").Value = 0.1148148148148148
$ws.Range("R6").Value = 0.05925925925925926
$ws.Range("S6").Value = 0.3444444444444444
$ws.Range("B7").Value = 0.1517509727626459
$ws.Range("D7").Value = 0.007782101167315175
$ws.Range("F7").Value = 0.007782101167315175
$ws.Range("J7").Value = 0.1478599221789883
$ws.Range("O7").Value = 0.01945525291828794
$ws.Range("Q7").Value = 0.1167315175097276
$ws.Range("R7").Value = 0.08560311284046693
$ws.Range("S7").Value = 0.4630350194552529
$ws.Range("B8").Value = 0.1292517006802721
$ws.Range("D8").Value = 0.0272108843537415
$ws.Range("E8").Value = 0.006802721088435374
$ws.Range("F8").Value = 0.04308390022675737
$ws.Range("J8").Value = 0.1292517006802721
$ws.Range("O8").Value = 0.01587301587301587
$ws.Range("Q8").Value = 0.1859410430839002
$ws.Range("R8").Value = 0.05442176870748299
$ws.Range("S8").Value = 0.4081632653061225
$ws.Range("B9").Value = 0.1494845360824742
$ws.Range("D9").Value = 0.01030927835051546
$ws.Range("F9").Value = 0.05154639175257732
$ws.Range("J9").Value = 0.1134020618556701
$ws.Range("O9").Value = 0.02577319587628866
$ws.Range("Q9").Value = 0.2061855670103093
$ws.Range("R9").Value = 0.1082474226804124
$ws.Range("S9").Value = 0.3350515463917526
$ws.Range("B10").Value = 0.1447631754503002
$ws.Range("D10").Value = 0.0200133422281521
$ws.Range("F10").Value = 0.0780520346897932
$ws.Range("J10").Value = 0.1034022681787859
$ws.Range("O10").Value = 0.01467645096731154
$ws.Range("Q10").Value = 0.1894596397598399
$ws.Range("R10").Value = 0.08005336891260841
$ws.Range("S10").Value = 0.3695797198132088
$ws.Range("G11").Value = 0.1462765957446809
$ws.Range("J11").Value = 0.07180851063829788
$ws.Range("K11").Value = 0.1888297872340426
$ws.Range("L11").Value = 0.5824468085106383
$ws.Range("S11").Value = 0.01063829787234043
$ws.Range("G12").Value = 0.7610619469026548
$ws.Range("J12").Value = 0.1902654867256637
$ws.Range("K12").Value = 0.004424778761061947
$ws.Range("L12").Value = 0.02654867256637168
$ws.Range("S12").Value = 0.01769911504424779
$ws.Range("G13").Value = 0.72
$ws.Range("J13").Value = 0.26
$ws.Range("S13").Value = 0.02
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
$ws.Range("F15").Value = 0.0158102766798419
$ws.Range("H15").Value = 0.1067193675889328
$ws.Range("I15").Value = 0.05928853754940711
$ws.Range("J15").Value = 0.4031620553359684
$ws.Range("K15").Value = 0.07905138339920949
$ws.Range("M15").Value = 0.007905138339920948
$ws.Range("N15").Value = 0.003952569169960474
$ws.Range("O15").Value = 0.09881422924901186
$ws.Range("S15").Value = 0.225296442687747
$ws.Range("F16").Value = 0.02127659574468085
$ws.Range("H16").Value = 0.1702127659574468
$ws.Range("I16").Value = 0.05673758865248227
$ws.Range("J16").Value = 0.4148936170212766
$ws.Range("K16").Value = 0.1099290780141844
$ws.Range("M16").Value = 0.01418439716312057
$ws.Range("N16").Value = 0.003546099290780142
$ws.Range("O16").Value = 0.06028368794326241
$ws.Range("S16").Value = 0.148936170212766
$ws.Range("F17").Value = 0.02991452991452992
$ws.Range("H17").Value = 0.1581196581196581
$ws.Range("I17").Value = 0.08974358974358974
$ws.Range("J17").Value = 0.4294871794871795
$ws.Range("K17").Value = 0.09829059829059829
$ws.Range("M17").Value = 0.02777777777777778
$ws.Range("N17").Value = 0.002136752136752137
$ws.Range("O17").Value = 0.04700854700854701
$ws.Range("S17").Value = 0.1175213675213675
$ws.Range("F18").Value = 0.02475247524752475
$ws.Range("H18").Value = 0.1584158415841584
$ws.Range("I18").Value = 0.07425742574257425
$ws.Range("J18").Value = 0.4653465346534654
$ws.Range("K18").Value = 0.103960396039604
$ws.Range("M18").Value = 0.01485148514851485
$ws.Range("N18").Value = 0.004950495049504951
$ws.Range("O18").Value = 0.06930693069306931
$ws.Range("S18").Value = 0.08415841584158416
$ws.Range("F19").Value = 0.02908587257617729
$ws.Range("H19").Value = 0.1842105263157895
$ws.Range("I19").Value = 0.07548476454293629
$ws.Range("J19").Value = 0.3725761772853186
$ws.Range("K19").Value = 0.1267313019390582
$ws.Range("M19").Value = 0.02008310249307479
$ws.Range("N19").Value = 0.001385041551246537
$ws.Range("O19").Value = 0.07409972299168975
$ws.Range("S19").Value = 0.1163434903047091
